{"js": "// Apply the CV edits described by the diff:\n//  1. Shorten the D&A practice bullet (drop \"direct $XXM+ portfolio with \").\n//  2. Delete the \"Designed credit risk AI models...\" bullet entirely.\n//  3. Rename \"Various Companies\" -> \"Microsoft, UTU & Others\".\n//  4. Rename \"Software Engineering & Technical Consulting\" -> \"...Technical Leadership\".\n//  5. Rewrite the \"Progressive roles...\" summary paragraph.\n//  6. Trim \" JAPAC Hackathon Winner.\" off the CatchMe blurb.\n\nasync function replaceText(context, findText, newText, options) {\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(findText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Trim the CSAT bullet.\nawait replaceText(\n  context,\n  \"Built D&A practice from 0 to 1 across 6 countries; direct $XXM+ portfolio with 97% CSAT\",\n  \"Built D&A practice from 0 to 1 across 6 countries; 97% CSAT\"\n);\n\n// 2. Delete the \"Designed credit risk AI models...\" bullet paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Designed credit risk AI models improving accuracy by 15% with alternative data sources\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the 'Designed credit risk AI models...' paragraph\");\n}\ntargetParagraph.delete();\nawait context.sync();\n\n// 3. \"Various Companies\" -> \"Microsoft, UTU & Others\"\nawait replaceText(context, \"Various Companies\", \"Microsoft, UTU & Others\");\n\n// 4. \"Software Engineering & Technical Consulting\" -> \"...Technical Leadership\"\nawait replaceText(\n  context,\n  \"Software Engineering & Technical Consulting\",\n  \"Software Engineering & Technical Leadership\"\n);\n\n// 5. Rewrite the \"Progressive roles...\" summary paragraph.\nawait replaceText(\n  context,\n  \"Progressive roles in software development, systems integration, and consulting in financial services and algorithmic trading.\",\n  \"Windows Kernel development (Microsoft), payment systems (UTU Singapore), founded Truckaurbus B2B marketplace.\"\n);\n\n// 6. Trim \" JAPAC Hackathon Winner.\" from the CatchMe blurb.\nawait replaceText(\n  context,\n  \"Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency. JAPAC Hackathon Winner.\",\n  \"Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency.\"\n);\n", "ps1": "# Apply the CV edits described by the diff:\n#  1. Shorten the D&A practice bullet (drop \"direct $XXM+ portfolio with \").\n#  2. Delete the \"Designed credit risk AI models...\" bullet entirely.\n#  3. Rename \"Various Companies\" -> \"Microsoft, UTU & Others\".\n#  4. Rename \"Software Engineering & Technical Consulting\" -> \"...Technical Leadership\".\n#  5. Rewrite the \"Progressive roles...\" summary paragraph.\n#  6. Trim \" JAPAC Hackathon Winner.\" off the CatchMe blurb.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphContaining {\n    param($doc, [string]$needle)\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Text -like \"*$needle*\") {\n            return $p\n        }\n    }\n    throw \"No paragraph found containing: $needle\"\n}\n\n# 1. Trim the CSAT bullet (whole paragraph is a single run).\n$p1 = Get-ParagraphContaining $d 'Built D&A practice from 0 to 1 across 6 countries'\n$p1.Range.Text = 'Built D&A practice from 0 to 1 across 6 countries; 97% CSAT'\n\n# 2. Delete the \"Designed credit risk AI models...\" bullet paragraph entirely.\n$p2 = Get-ParagraphContaining $d 'Designed credit risk AI models improving accuracy by 15% with alternative data sources'\n$p2.Range.Delete()\n\n# 3. \"Various Companies\" -> \"Microsoft, UTU & Others\" (whole paragraph).\n$p3 = Get-ParagraphContaining $d 'Various Companies'\n$p3.Range.Text = 'Microsoft, UTU & Others'\n\n# 4. \"Software Engineering & Technical Consulting\" -> \"...Technical Leadership\".\n#    This paragraph also holds the separate \"  \" / italic-date runs, so only\n#    replace the leading run's text via a Find scoped to the paragraph range.\n$p4 = Get-ParagraphContaining $d 'Software Engineering & Technical Consulting'\n$p4Range = $p4.Range\n$p4Find = $p4Range.Find\n$p4Find.ClearFormatting()\n$p4Found = $p4Find.Execute('Software Engineering & Technical Consulting', $false, $false, $false, $false, $false, $true, 0, $false)\nif (-not $p4Found) {\n    throw \"Could not locate the 'Software Engineering & Technical Consulting' run\"\n}\n$p4Range.Text = 'Software Engineering & Technical Leadership'\n\n# 5. Rewrite the \"Progressive roles...\" summary paragraph (single run).\n$p5 = Get-ParagraphContaining $d 'Progressive roles in software development'\n$p5.Range.Text = 'Windows Kernel development (Microsoft), payment systems (UTU Singapore), founded Truckaurbus B2B marketplace.'\n\n# 6. Trim \" JAPAC Hackathon Winner.\" from the CatchMe blurb (single run).\n$p6 = Get-ParagraphContaining $d 'Agentic AI for enterprise trust decisions'\n$p6.Range.Text = 'Agentic AI for enterprise trust decisions. APLS pattern learning + cascade routing. 86% cost reduction, sub-50ms latency.'\n\nWrite-Output \"done\"\n"}
